$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.487.61"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "1.879.72"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "0.7168"
$ws.Range("E5").Value = "  +1.97%  "
$ws.Range("D6").Value = "241.87"
$ws.Range("E6").Value = "  +1.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9990"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "0.07899"
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("E9").Value = "  +3.49%  "
$ws.Range("D10").Value = "25.35"
$ws.Range("E10").Value = "  +7.57%  "
$ws.Range("D11").Value = "0.08275"
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "0.7332"
$ws.Range("E12").Value = "  +4.00%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.886.17"
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("D14").Value = "5.301"
$ws.Range("E14").Value = "  +2.21%  "
$ws.Range("D15").Value = "91.48"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").Value = "29.503.43"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.960"
$ws.Range("E17").Value = "  +2.42%  "
$ws.Range("D18").Value = "247.87"
$ws.Range("E18").Value = "  +4.85%  "
$ws.Range("D19").Value = "0.000007874"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").Value = "13.38"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "8.026"
$ws.Range("E21").Value = "  +7.05%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "0.9987"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "0.9988"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  +13.81%  "
$ws.Range("D25").Value = "163.82"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("D26").Value = "9.069"
$ws.Range("E26").Value = "  +2.34%  "
$ws.Range("D27").Value = "18.37"
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("E28").Value = "  -2.76%  "
$ws.Range("D29").Value = "1.504"
$ws.Range("E29").Value = "  +2.29%  "
$ws.Range("D30").Value = "4.395"
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("D31").Value = "4.125"
$ws.Range("E31").Value = "  +2.79%  "
$ws.Range("D32").Value = "0.05291"
$ws.Range("E32").Value = "  +2.46%  "
$ws.Range("E33").Value = "  +2.42%  "
$ws.Range("E34").Value = "  +3.05%  "
$ws.Range("D35").Value = "0.7299"
$ws.Range("E35").Value = "  +2.78%  "
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "0.01873"
$ws.Range("E37").Value = "  +1.43%  "
$ws.Range("D38").Value = "1.227.70"
$ws.Range("E38").Value = "  +6.11%  "
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "0.9117"
$ws.Range("E40").Value = "  -2.05%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "74.88"
$ws.Range("E41").Value = "  +6.69%  "
$ws.Range("D42").Value = "6.203"
$ws.Range("E42").Value = "  +3.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9990"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "102.82"
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("D45").Value = "2.047.25"
$ws.Range("E45").Value = "  +3.17%  "
$ws.Range("D46").Value = "0.5259"
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("D47").Value = "2.964"
$ws.Range("E47").Value = "  +14.18%  "
$ws.Range("D48").Value = "1.783"
$ws.Range("E48").Value = "  +2.51%  "
$ws.Range("D49").Value = "9.353"
$ws.Range("E49").Value = "  +2.27%  "
$ws.Range("D50").Value = "0.4344"
$ws.Range("E50").Value = "  +2.15%  "
$ws.Range("E51").Value = "  +2.28%  "
